# Apply the MeARM Inverse Kinematics worksheet update:
#  - Update measured/supplied angle data (columns B and J) for rows 15-18, 20-23,
#    25-28, 30-33 on Sheet1, and recompute the dependent formula columns
#    (C, D, E, F, K, L, M, N) to match the new calibration numbers.
#  - Clear the leftover placeholder notes in D15/G15 ("96(not acc)" / "impossible").
#  - Update the sheet view (zoom level and selected cell).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New "Supplied" values for the Left servo (B) and Right servo (J) columns.
$rowData = @{
  15 = @{ B = 80;  J = 45 }
  16 = @{ B = 65;  J = 60 }
  17 = @{ B = 50;  J = 75 }
  18 = @{ B = 35;  J = 90 }
  20 = @{ B = 95;  J = 40 }
  21 = @{ B = 73;  J = 57 }
  22 = @{ B = 60;  J = 70 }
  23 = @{ B = 50;  J = 85 }
  25 = @{ B = 105; J = 40 }
  26 = @{ B = 85;  J = 55 }
  27 = @{ B = 75;  J = 70 }
  28 = @{ B = 65;  J = 85 }
  30 = @{ B = 120; J = 40 }
  31 = @{ B = 103; J = 56 }
  32 = @{ B = 87;  J = 70 }
  33 = @{ B = 77;  J = 85 }
}

foreach ($r in $rowData.Keys) {
  $b = $rowData[$r].B
  $j = $rowData[$r].J

  $ws.Range("B$r").Value = $b
  $ws.Range("J$r").Value = $j

  $ws.Range("C$r").Formula = "=A$r-B$r"
  $ws.Range("D$r").Formula = "=(81/SIN(RADIANS(L$r)))*SIN(RADIANS(B$r))"
  $ws.Range("E$r").Formula = "=J$r"
  $ws.Range("F$r").Formula = "=180-E$r-B$r"

  $ws.Range("K$r").Formula = "=I$r-J$r"
  $ws.Range("L$r").Formula = "=(180-B$r)/2"
  $ws.Range("N$r").Formula = "=90-J$r"
  $ws.Range("M$r").Formula = "=90-L$r-N$r"
}

# Remove the stray "impossible" note left in G15 (D15's placeholder text
# was already overwritten by the formula set in the loop above).
$ws.Range("G15").Value = ""

# The separator rows (19, 24, 29) no longer carry the placeholder formulas.
foreach ($r in 19, 24, 29) {
  foreach ($col in "C", "D", "E", "F", "K", "L", "M", "N") {
    $ws.Range("$col$r").Value = ""
  }
}

# Update the view: drop the old scroll position, change zoom, move the
# active selection to H40.
$ws.Activate()
$win = $excel.ActiveWindow
$win.Zoom = 88
$ws.Range("H40").Select()
